# Trade #3 closed at 2026-02-16 21:50:53 - leadlag UP +0.000%
# A new trade-log row (row 4) is appended to both the "All Trades" and the
# "leadlag" sheets with identical data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "leadlag")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(4, 1).Value = 3

    # Force column B to stay text so the date-like string "2026-02-16"
    # isn't auto-converted into a date serial number, then drop the
    # temporary number format again so the cell keeps the sheet's default
    # (unstyled) look, matching the other rows.
    $ws.Cells.Item(4, 2).NumberFormat = "@"
    $ws.Cells.Item(4, 2).Value = "2026-02-16"
    $ws.Cells.Item(4, 2).ClearFormats()

    $ws.Cells.Item(4, 3).Value = "21:50:53"
    $ws.Cells.Item(4, 4).Value = "leadlag"
    $ws.Cells.Item(4, 5).Value = "UP"
    $ws.Cells.Item(4, 6).Value = 68369.25999999999

    # Exit Price is blank (trade still open) - keep it an empty text cell.
    $ws.Cells.Item(4, 7).Formula = '=""'

    $ws.Cells.Item(4, 8).Value = "OPEN"
    $ws.Cells.Item(4, 9).Value = 0
    $ws.Cells.Item(4, 10).Value = 0
    $ws.Cells.Item(4, 11).Value = 100
    $ws.Cells.Item(4, 12).Value = 0.75
    $ws.Cells.Item(4, 13).Value = "Coinbase leading with 0.109% move"

    # Exit Reason is blank (trade still open) - keep it an empty text cell.
    $ws.Cells.Item(4, 14).Formula = '=""'

    $ws.Cells.Item(4, 15).Value = 0
}
